# Commit: Add choice to practise (verb forms or words) - unique selection
# Fixes verb conjugation typos, adds gendered articles to existing vocab,
# and appends a large new vocabulary section (bedroom/bathroom/living room/
# kitchen/garage/prepositions/misc) to the 'words' sheet.

$wb = $excel.ActiveWorkbook
$wsVerbs = $wb.Worksheets.Item("verbs")
$wsWords = $wb.Worksheets.Item("words")

# --- 1) Fix typo'd conjugations of 'mostrar' (to show) on the verbs sheet, row 48 ---
$wsVerbs.Range("B48").Value = "muestro"
$wsVerbs.Range("C48").Value = "mostrás"
$wsVerbs.Range("D48").Value = "muestra"
$wsVerbs.Range("E48").Value = "mostramos"
$wsVerbs.Range("F48").Value = "muestran"

# --- 2) Add missing gendered articles (el/la/las) to the existing classroom-object words ---
$wsWords.Range("A1").Value = 'el cuaderno'
$wsWords.Range("A2").Value = 'la regla'
$wsWords.Range("A3").Value = 'la cinta de papel'
$wsWords.Range("A4").Value = 'la calculadora'
$wsWords.Range("A5").Value = 'el lápiz'
$wsWords.Range("A6").Value = 'la plasticola'
$wsWords.Range("A7").Value = 'el pincel'
$wsWords.Range("A8").Value = 'el lapicero'
$wsWords.Range("A9").Value = 'la tijera'
$wsWords.Range("A10").Value = 'el libro'
$wsWords.Range("A11").Value = 'el resaltador'
$wsWords.Range("A12").Value = 'la perforadora'
$wsWords.Range("A13").Value = 'la goma'
$wsWords.Range("A14").Value = 'las hojas de papel'
$wsWords.Range("A15").Value = 'el sacapuntas'
$wsWords.Range("A16").Value = 'el escritorio'
$wsWords.Range("A17").Value = 'el tacho de basura'
$wsWords.Range("A18").Value = 'el reloj'
$wsWords.Range("A19").Value = 'la abrochadora'
$wsWords.Range("A20").Value = 'la carpeta'
$wsWords.Range("A21").Value = 'la computadora portátil'
$wsWords.Range("A22").Value = 'la agenda'
$wsWords.Range("A23").Value = 'la lapicera'
$wsWords.Range("A24").Value = 'la impresora'
$wsWords.Range("A25").Value = 'el bolígrafo'

# --- 3) Append the new vocabulary rows (bedroom, bathroom, living room, kitchen,
#        garage, prepositions of place, and assorted extra words) starting row 26 ---
$newRows = @(
    @('el dormitorio', 'bedroom'),
    @('la cama', 'bed'),
    @('la lámpara', 'lamp'),
    @('la mesita de luz', 'nightstand'),
    @('la cuadro', 'picture/painting[c]'),
    @('la pintura', 'painting [p]'),
    @('el espejo', 'mirror'),
    @('el quarto de baño', 'bathroom'),
    @('la bañera', 'bathtub'),
    @('el inodoro', 'toilet'),
    @('la alfombrilla', 'mat'),
    @('el lavabo', 'sink (to wash hands)'),
    @('la ducha', 'shower'),
    @('la toalla', 'towel'),
    @('el salón', 'living room [s]'),
    @('la puerta', 'door'),
    @('el sillón', 'arnchair'),
    @('el sofa', 'sofa'),
    @('el suelo', 'soil'),
    @('el piso', 'floor'),
    @('la silla', 'chair'),
    @('la planta', 'plant'),
    @('la mesa', 'table'),
    @('el armario', 'cabinet'),
    @('la estantería', 'shelf'),
    @('el techo', 'ceiling/roof'),
    @('la chimenea', 'chimney'),
    @('la ventana', 'window'),
    @('la cocina', 'kitchen'),
    @('la heladera', 'refrigerator'),
    @('el horno', 'oven'),
    @('el lavarropas', 'washing machine'),
    @('el microondas', 'microwave'),
    @('el lavatorio', 'sink (for kitchen)'),
    @('el lavavajillas', 'dishwasher'),
    @('el cajón', 'drawer'),
    @('el auto', 'car'),
    @('el garaje', 'garage'),
    @('el baño', 'bath'),
    @('el living', 'living room [l]'),
    @('en el centro de', 'in the center of'),
    @('encima de', 'above of/on [e]'),
    @('sobre', 'on [s]'),
    @('abajo de', 'under of'),
    @('a la derecha de', 'to the right of'),
    @('a la izquierda de', 'to the left of'),
    @('entre y', 'between and'),
    @('al lado de', 'next to'),
    @('adelante de', 'in front of'),
    @('atrás de', 'behind of'),
    @('alrededor de', 'around (the house)'),
    @('la bicicleta', 'bicycle'),
    @('la pelota', 'ball'),
    @('las zapatos', 'shoes'),
    @('el ratón', 'mouse'),
    @('la serpiente', 'snake'),
    @('el pájaro', 'bird'),
    @('la falda', 'skirt'),
    @('la banda', '(music) band'),
    @('la maleta', 'suitcase'),
    @('el periódico', 'newspaper'),
    @('el armario', 'closet'),
    @('el bocadillo', 'snack'),
    @('la haladera', 'refrigerator'),
    @('la libreta', 'notepad'),
    @('la camisa', 'shirt')
)

$startRow = 26
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $wsWords.Cells.Item($r, 1).Value = $newRows[$i][0]
    $wsWords.Cells.Item($r, 2).Value = $newRows[$i][1]
}

# --- 4) Widen column A on 'words' to fit the longer (articled) entries ---
$wsWords.Columns.Item(1).ColumnWidth = 21.42

# --- 5) Restore view state: scroll/selection on verbs, then make 'words' the active tab again ---
$wsVerbs.Range("F49").Select()
$wsWords.Activate()
$wsWords.Range("B92").Select()

